$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text (string) type, since the source
# workbook stores these as inline strings (e.g. "1.00", "214.45") rather
# than numbers. Force text format before assigning so Excel does not
# auto-convert/round them to numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.701.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.622.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.601.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.703.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0729"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.450.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.573"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.951"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.763.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.764"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Algorand"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0965"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.56%  "
